$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "260.74"),
    @("E2", "1.92%"),
    @("D3", "27.24"),
    @("E3", "3.14%"),
    @("D4", "4.687"),
    @("E4", "0.82%"),
    @("D5", "0.06124"),
    @("E5", "3.38%"),
    @("D6", "6.659"),
    @("E6", "0.67%"),
    @("E7", "-0.08%"),
    @("D8", "0.9210"),
    @("E8", "1.33%"),
    @("D9", "0.1401"),
    @("E9", "1.64%"),
    @("D10", "0.04717"),
    @("E10", "14.54%"),
    @("D11", "0.07090"),
    @("E11", "1.37%"),
    @("D12", "0.03053"),
    @("E12", "0.42%"),
    @("D13", "0.09061"),
    @("E13", "-0.27%"),
    @("D14", "0.001529"),
    @("E14", "-0.21%"),
    @("D15", "0.0006080"),
    @("E15", "0.85%"),
    @("D16", "0.006045"),
    @("E16", "-0.12%"),
    @("E17", "-0.56%"),
    @("E18", "-0.01%"),
    @("D19", "2.163"),
    @("E19", "-0.61%"),
    @("D21", "0.1304"),
    @("E21", "1.52%"),
    @("D22", "4.098"),
    @("E22", "6.29%"),
    @("D23", "0.04240"),
    @("E23", "0.51%"),
    @("D24", "0.001216"),
    @("E24", "0.14%"),
    @("D25", "0.003801"),
    @("E25", "-18.87%"),
    @("E26", "0.09%"),
    @("D27", "0.0001574"),
    @("E27", "3.35%"),
    @("D40", "0.03865"),
    @("E40", "2.36%"),
    @("D41", "0.1113"),
    @("E41", "1.77%"),
    @("D42", "0.004077"),
    @("E42", "9.10%"),
    @("D43", "0.01633"),
    @("E43", "13.04%"),
    @("D44", "0.002216"),
    @("E44", "-9.08%"),
    @("D45", "0.00005154"),
    @("E45", "0.02%"),
    @("E46", "-0.01%"),
    @("E47", "8.03%"),
    @("D48", "0.1355"),
    @("E48", "-43.74%"),
    @("E49", "-0.01%"),
    @("E50", "-0.01%")
)

foreach ($item in $changes) {
    $cellref = $item[0]
    $val = $item[1]
    $ws.Range($cellref).NumberFormat = "@"
    $ws.Range($cellref).Value = $val
    $ws.Range($cellref).ClearFormats()
}

